$wb = $excel.ActiveWorkbook

# --- 1_Vocab_Ex sheet (sheet2 / index 2): move selection, no longer the active tab ---
$wsVocabEx = $wb.Worksheets.Item(2)
$wsVocabEx.Activate()
$wsVocabEx.Range("C2:C5").Select()

# --- 2_ -> 2_Vocab sheet (sheet3 / index 3): rename, populate data, resize column, activate ---
$ws2Vocab = $wb.Worksheets.Item(3)
$ws2Vocab.Name = "2_Vocab"

# NOTE: shared-string order in the source workbook shows the table body (A2:A7)
# was entered before the header question (A1), so we replay writes in that
# same order to reproduce identical sharedStrings.xml indices.
$ws2Vocab.Range("A2").Value = "The rate of arrival of new cars"
$ws2Vocab.Range("B2").Value = "C"
$ws2Vocab.Range("C2").Value = "Independent Variable"

$ws2Vocab.Range("A3").Value = "The length of red lights"
$ws2Vocab.Range("B3").Value = "A"
$ws2Vocab.Range("C3").Value = "State Variable"

$ws2Vocab.Range("A4").Value = "The use of turn arrows"
$ws2Vocab.Range("B4").Value = "A"
$ws2Vocab.Range("C4").Value = "Parameter"

$ws2Vocab.Range("A5").Value = "The average wait time for each car"
$ws2Vocab.Range("B5").Value = "D"
$ws2Vocab.Range("C5").Value = "Metric"

$ws2Vocab.Range("A6").Value = "The number of cars in a given lane at 5:14 pm"
$ws2Vocab.Range("B6").Value = "B"

$ws2Vocab.Range("A7").Value = "The maximum wait time for a driver"
$ws2Vocab.Range("B7").Value = "D"

$ws2Vocab.Range("A1").Value = "Imagine you are modeling traffic flow at a busy corner.    You want as little back up of traffic as possible during evening rush hour.  There is no room to increase the number of lanes, but you can alter the pattern of the lights.  What type of variable is each of these?"

$ws2Vocab.Rows.Item(1).RowHeight = 120
$ws2Vocab.Rows.Item(6).RowHeight = 30
$ws2Vocab.Rows.Item(7).RowHeight = 30

$ws2Vocab.Columns.Item(1).ColumnWidth = 31.8

$ws2Vocab.Activate()
